# Re-apply the latest scraped crypto price/volume snapshot (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force plain numeric-looking strings (e.g. "1.000",
# "0.7008") to land as literal text instead of being auto-converted to
# numbers by Value-assignment type inference. A leading apostrophe marks
# the scratch value as text; copy/PasteSpecial(values-only) then carries
# just that literal text into the destination cell, leaving the
# destination cells own formatting untouched. The scratch column is
# deleted again at the end so no extra column is left behind.
$scratch = $ws.Range("ZZ1")

$ws.Range("D2").Value = '29.315.81'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '1.860.43'
$ws.Range("E3").Value = '  +0.11%  '
$scratch.Value = "'1.000"
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = '  +0.04%  '
$scratch.Value = "'0.7008"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -0.17%  '
$scratch.Value = "'237.60"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E7").Value = '  +0.11%  '
$scratch.Value = "'0.07781"
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  -3.24%  '
$scratch.Value = "'0.3040"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  -0.16%  '
$scratch.Value = "'24.71"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +6.18%  '
$scratch.Value = "'0.08149"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '1.858.02'
$ws.Range("E12").Value = '  -0.17%  '
$scratch.Value = "'5.206"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  +0.43%  '
$scratch.Value = "'0.7128"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = '  -0.39%  '
$scratch.Value = "'89.19"
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("D16").Value = '29.294.03'
$ws.Range("E16").Value = '  +0.11%  '
$scratch.Value = "'5.775"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +0.18%  '
$scratch.Value = "'242.40"
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +2.64%  '
$scratch.Value = "'0.000007762"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  -0.91%  '
$scratch.Value = "'13.15"
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '2.094.07'
$ws.Range("E22").Value = '  -0.69%  '
$ws.Range("E23").Value = '  +0.09%  '
$scratch.Value = "'7.510"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +0.75%  '
$scratch.Value = "'162.09"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +0.13%  '
$scratch.Value = "'8.857"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -1.32%  '
$scratch.Value = "'0.1432"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -1.48%  '
$scratch.Value = "'18.05"
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -0.02%  '
$scratch.Value = "'1.896"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  -4.90%  '
$scratch.Value = "'1.367"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -4.65%  '
$scratch.Value = "'1.473"
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -0.85%  '
$scratch.Value = "'4.288"
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  -2.65%  '
$scratch.Value = "'4.022"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -0.76%  '
$scratch.Value = "'0.05156"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  -1.05%  '
$scratch.Value = "'1.179"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  +0.74%  '
$scratch.Value = "'0.7032"
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  -0.77%  '
$scratch.Value = "'0.9936"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  -0.87%  '
$scratch.Value = "'2.678"
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +0.61%  '
$scratch.Value = "'0.01844"
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -0.17%  '
$scratch.Value = "'2.690"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("D41").Value = '1.174.26'
$ws.Range("E41").Value = '  +2.30%  '
$scratch.Value = "'0.9122"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -1.45%  '
$scratch.Value = "'6.012"
$scratch.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +0.93%  '
$scratch.Value = "'71.36"
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +0.65%  '
$scratch.Value = "'0.4230"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -1.04%  '
$scratch.Value = "'1.001"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +0.17%  '
$scratch.Value = "'101.18"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  -1.98%  '
$scratch.Value = "'0.5350"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -1.21%  '
$scratch.Value = "'1.743"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  -2.05%  '
$scratch.Value = "'9.115"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  -0.77%  '
$scratch.Value = "'6.919"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  -0.26%  '

# Clean up the scratch column entirely (also clears clipboard marching-ants state).
$scratch.Value = $null
$ws.Columns.Item($scratch.Column).Delete()
